$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.Value = "'" + $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '42.488.33'
Set-TextValue $ws.Range("E2") '  -2.59%  '

Set-TextValue $ws.Range("D3") '2.278.31'
Set-TextValue $ws.Range("E3") '  -4.13%  '

Set-TextValue $ws.Range("E4") '  +0.24%  '

Set-TextValue $ws.Range("D5") '300.69'
Set-TextValue $ws.Range("E5") '  -2.90%  '

Set-TextValue $ws.Range("D6") '97.23'
Set-TextValue $ws.Range("E6") '  -7.12%  '

Set-TextValue $ws.Range("D7") '0.504'
Set-TextValue $ws.Range("E7") '  -0.65%  '

Set-TextValue $ws.Range("E8") '  +0.23%  '

Set-TextValue $ws.Range("D9") '0.500'
Set-TextValue $ws.Range("E9") '  -3.92%  '

Set-TextValue $ws.Range("D10") '34.12'
Set-TextValue $ws.Range("E10") '  -5.31%  '

Set-TextValue $ws.Range("D11") '0.0788'
Set-TextValue $ws.Range("E11") '  -3.06%  '

Set-TextValue $ws.Range("D12") '50.65'
Set-TextValue $ws.Range("E12") '  -5.19%  '

Set-TextValue $ws.Range("E13") '  +0.27%  '

Set-TextValue $ws.Range("D14") '6.68'
Set-TextValue $ws.Range("E14") '  -4.43%  '

Set-TextValue $ws.Range("D15") '2.644.71'
Set-TextValue $ws.Range("E15") '  -3.50%  '

Set-TextValue $ws.Range("D16") '15.29'
Set-TextValue $ws.Range("E16") '  -2.25%  '

Set-TextValue $ws.Range("D17") '2.281.15'
Set-TextValue $ws.Range("E17") '  -3.93%  '

Set-TextValue $ws.Range("D18") '0.787'
Set-TextValue $ws.Range("E18") '  -2.88%  '

Set-TextValue $ws.Range("D19") '42.481.76'
Set-TextValue $ws.Range("E19") '  -2.49%  '

Set-TextValue $ws.Range("D20") '11.54'
Set-TextValue $ws.Range("E20") '  -3.00%  '

Set-TextValue $ws.Range("D21") '0.0₃0894'
Set-TextValue $ws.Range("E21") '  -2.58%  '

Set-TextValue $ws.Range("D22") '5.99'
Set-TextValue $ws.Range("E22") '  -5.16%  '

Set-TextValue $ws.Range("D23") '66.68'
Set-TextValue $ws.Range("E23") '  -2.56%  '

Set-TextValue $ws.Range("D24") '234.10'
Set-TextValue $ws.Range("E24") '  -2.85%  '

Set-TextValue $ws.Range("D25") '1.94'
Set-TextValue $ws.Range("E25") '  -5.31%  '

Set-TextValue $ws.Range("D26") '2.49'
Set-TextValue $ws.Range("E26") '  -4.75%  '

Set-TextValue $ws.Range("E27") '  -0.14%  '

Set-TextValue $ws.Range("D28") '24.51'
Set-TextValue $ws.Range("E28") '  -5.02%  '

Set-TextValue $ws.Range("D29") '2.18'
Set-TextValue $ws.Range("E29") '  +3.29%  '

Set-TextValue $ws.Range("D30") '34.02'
Set-TextValue $ws.Range("E30") '  -7.06%  '

Set-TextValue $ws.Range("D31") '165.03'
Set-TextValue $ws.Range("E31") '  +2.59%  '

Set-TextValue $ws.Range("D32") '9.08'
Set-TextValue $ws.Range("E32") '  -4.79%  '

Set-TextValue $ws.Range("E33") '  +0.28%  '

Set-TextValue $ws.Range("D34") '4.95'
Set-TextValue $ws.Range("E34") '  -5.86%  '

Set-TextValue $ws.Range("E35") '  -4.76%  '

Set-TextValue $ws.Range("D36") '0.0697'
Set-TextValue $ws.Range("E36") '  -5.82%  '

Set-TextValue $ws.Range("D37") '4.38'
Set-TextValue $ws.Range("E37") '  -5.81%  '

Set-TextValue $ws.Range("D38") '2.82'
Set-TextValue $ws.Range("E38") '  -9.59%  '

Set-TextValue $ws.Range("D39") '16.13'
Set-TextValue $ws.Range("E39") '  -11.80%  '

Set-TextValue $ws.Range("B40") 'Kaspa'
Set-TextValue $ws.Range("C40") 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D40") '0.100'
Set-TextValue $ws.Range("E40") '  -5.03%  '

Set-TextValue $ws.Range("B41") 'ARBITRUM'
Set-TextValue $ws.Range("C41") 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D41") '1.77'
Set-TextValue $ws.Range("E41") '  -8.63%  '

Set-TextValue $ws.Range("E42") '  -3.69%  '

Set-TextValue $ws.Range("D43") '2.39'
Set-TextValue $ws.Range("E43") '  -8.04%  '

Set-TextValue $ws.Range("D44") '1.965.91'
Set-TextValue $ws.Range("E44") '  -3.29%  '

Set-TextValue $ws.Range("D45") '0.0282'
Set-TextValue $ws.Range("E45") '  -2.79%  '

Set-TextValue $ws.Range("D46") '17.91'
Set-TextValue $ws.Range("E46") '  -9.35%  '

Set-TextValue $ws.Range("D47") '9.73'
Set-TextValue $ws.Range("E47") '  -8.00%  '

Set-TextValue $ws.Range("D48") '2.83'
Set-TextValue $ws.Range("E48") '  -9.49%  '

Set-TextValue $ws.Range("B49") 'HuobiToken'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range("D49") '2.85'
Set-TextValue $ws.Range("E49") '  -4.07%  '

Set-TextValue $ws.Range("B50") 'THORChain'
Set-TextValue $ws.Range("C50") 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range("D50") '4.70'
Set-TextValue $ws.Range("E50") '  -0.80%  '

Set-TextValue $ws.Range("D51") '2.510.91'
Set-TextValue $ws.Range("E51") '  -3.64%  '
